# Update the "Table Caption" and "Image Caption" paragraph styles so that
# their captions render in bold, non-italic Times New Roman (matching the
# rest of the report body), and give the Image Caption style zero space
# after the paragraph (tightening up the figure/caption spacing), as part
# of the officedown-driven formatting refresh for tables and figures.

$d = $word.ActiveDocument

# --- "Table Caption" style (based on "Caption") ---------------------------
$tableCaption = $d.Styles("TableCaption")
$tableCaption.Font.Name = "Times New Roman"
$tableCaption.Font.Bold = $true
$tableCaption.Font.Italic = $false

# --- "Image Caption" style (based on "Caption") ----------------------------
$imageCaption = $d.Styles("ImageCaption")
$imageCaption.ParagraphFormat.SpaceAfter = 0
$imageCaption.Font.Name = "Times New Roman"
$imageCaption.Font.Bold = $true
$imageCaption.Font.Italic = $false
